$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 6634.1177
$ws.Range("J100").Value = 6768.75
$ws.Range("L100").Value = 6768.75
$ws.Range("N100").Value = -7850.75
$ws.Range("H107").Value = 649.55554
$ws.Range("I107").Value = 788.05
$ws.Range("J107").Value = 253.85715
$ws.Range("K107").Value = 788.05
$ws.Range("L107").Value = 253.85715
$ws.Range("M107").Value = 1131.95
$ws.Range("N107").Value = -4093.85715
$ws.Range("H129").Value = 1865
$ws.Range("J129").Value = 2625
$ws.Range("L129").Value = 7875
$ws.Range("N129").Value = -17875
$ws.Range("H132").Value = 3356.7678
$ws.Range("I132").Value = 1594.9556
$ws.Range("J132").Value = 10564.182
$ws.Range("K132").Value = 4784.8668
$ws.Range("L132").Value = 31692.546
$ws.Range("M132").Value = -2254.8668
$ws.Range("N132").Value = -36752.546
$ws.Range("H137").Value = 3204.8572
$ws.Range("I137").Value = 3375.1538
$ws.Range("J137").Value = 2928.125
$ws.Range("K137").Value = 10125.4614
$ws.Range("L137").Value = 8784.375
$ws.Range("M137").Value = -7575.4614
$ws.Range("N137").Value = -13884.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 60000
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61498
$ws.Range("H72").Value = 60000
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -187488
$ws.Range("H74").Value = 2052
$ws.Range("I74").Value = 1071.0588
$ws.Range("J74").Value = 4617.5386
$ws.Range("K74").Value = 1071.0588
$ws.Range("L74").Value = 4617.5386
$ws.Range("M74").Value = -197.0588
$ws.Range("N74").Value = -6365.5386
$ws.Range("H77").Value = 2052
$ws.Range("I77").Value = 1071.0588
$ws.Range("J77").Value = 4617.5386
$ws.Range("K77").Value = 5355.294
$ws.Range("L77").Value = 23087.693
$ws.Range("M77").Value = -987.2939999999999
$ws.Range("N77").Value = -31823.693
$ws.Range("H110").Value = 2318.5
$ws.Range("I110").Value = 2318.5
$ws.Range("K110").Value = 2318.5
$ws.Range("M110").Value = -273.5
$ws.Range("H122").Value = 4081.1428
$ws.Range("I122").Value = 3845.9443
$ws.Range("J122").Value = 5492.3335
$ws.Range("K122").Value = 11537.8329
$ws.Range("L122").Value = 16477.0005
$ws.Range("M122").Value = -9087.832900000001
$ws.Range("N122").Value = -21377.0005
$ws.Range("H132").Value = 13332.124
$ws.Range("I132").Value = 16002.543
$ws.Range("J132").Value = 3493.7368
$ws.Range("K132").Value = 48007.629
$ws.Range("L132").Value = 10481.2104
$ws.Range("M132").Value = -45477.629
$ws.Range("N132").Value = -15541.2104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4695.294
$ws.Range("I99").Value = 5272.857
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 5272.857
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -3774.857
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 1585.9744
$ws.Range("I105").Value = 1348.963
$ws.Range("J105").Value = 2119.25
$ws.Range("K105").Value = 1348.963
$ws.Range("L105").Value = 2119.25
$ws.Range("M105").Value = 398.037
$ws.Range("N105").Value = -5613.25
$ws.Range("H107").Value = 3392.3333
$ws.Range("I107").Value = 2869.8333
$ws.Range("J107").Value = 3914.8333
$ws.Range("K107").Value = 2869.8333
$ws.Range("L107").Value = 3914.8333
$ws.Range("M107").Value = -949.8332999999998
$ws.Range("N107").Value = -7754.8333
$ws.Range("H134").Value = 3832.743
$ws.Range("I134").Value = 3530.8147
$ws.Range("J134").Value = 4851.75
$ws.Range("K134").Value = 10592.4441
$ws.Range("L134").Value = 14555.25
$ws.Range("M134").Value = -8057.444100000001
$ws.Range("N134").Value = -19625.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3078.1836
$ws.Range("I31").Value = 2230.3235
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2230.3235
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1935.3235
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 3078.1836
$ws.Range("I34").Value = 2230.3235
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2230.3235
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2028.3235
$ws.Range("N34").Value = -5404
$ws.Range("H107").Value = 372.55
$ws.Range("I107").Value = 196.36363
$ws.Range("J107").Value = 587.8889
$ws.Range("K107").Value = 196.36363
$ws.Range("L107").Value = 587.8889
$ws.Range("M107").Value = 1723.63637
$ws.Range("N107").Value = -4427.8889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 14.4
$ws.Range("I12").Value = 18.5
$ws.Range("J12").Value = 10.615385
$ws.Range("K12").Value = 55.5
$ws.Range("L12").Value = 31.846155
$ws.Range("M12").Value = 117.5
$ws.Range("N12").Value = -377.846155
$ws.Range("H58").Value = 2738
$ws.Range("J58").Value = 3172.5
$ws.Range("L58").Value = 9517.5
$ws.Range("N58").Value = -9773.5
$ws.Range("H74").Value = 2833.3333
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -14122
$ws.Range("H77").Value = 2833.3333
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -46608
$ws.Range("H98").Value = 3199.6667
$ws.Range("I98").Value = 207.83333
$ws.Range("J98").Value = 6191.5
$ws.Range("K98").Value = 623.49999
$ws.Range("L98").Value = 18574.5
$ws.Range("M98").Value = 874.50001
$ws.Range("N98").Value = -21570.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 23814478
$ws.Range("J107").Value = 62500304
$ws.Range("L107").Value = 62500304
$ws.Range("N107").Value = -62504144
$ws.Range("H113").Value = 1396.1875
$ws.Range("I113").Value = 807.1429000000001
$ws.Range("J113").Value = 1854.3334
$ws.Range("K113").Value = 807.1429000000001
$ws.Range("L113").Value = 1854.3334
$ws.Range("M113").Value = 1362.8571
$ws.Range("N113").Value = -6194.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 665.35297
$ws.Range("I22").Value = 775.8570999999999
$ws.Range("J22").Value = 588
$ws.Range("K22").Value = 775.8570999999999
$ws.Range("L22").Value = 588
$ws.Range("M22").Value = -480.8570999999999
$ws.Range("N22").Value = -1178
$ws.Range("H27").Value = 665.35297
$ws.Range("I27").Value = 775.8570999999999
$ws.Range("J27").Value = 588
$ws.Range("K27").Value = 775.8570999999999
$ws.Range("L27").Value = 588
$ws.Range("M27").Value = -668.8570999999999
$ws.Range("N27").Value = -802
$ws.Range("H40").Value = 2367.25
$ws.Range("I40").Value = 1926.6
$ws.Range("J40").Value = 3101.6667
$ws.Range("K40").Value = 1926.6
$ws.Range("L40").Value = 3101.6667
$ws.Range("M40").Value = -1790.6
$ws.Range("N40").Value = -3373.6667
$ws.Range("H122").Value = 2163.077
$ws.Range("I122").Value = 2120
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 6360
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -3910
$ws.Range("N122").Value = -12100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1080
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 1991
$ws.Range("I132").Value = 1352.6923
$ws.Range("J132").Value = 2781.2856
$ws.Range("K132").Value = 4058.0769
$ws.Range("L132").Value = 8343.856800000001
$ws.Range("M132").Value = -1528.0769
$ws.Range("N132").Value = -13403.8568
$ws.Range("H136").Value = 12988623
$ws.Range("I136").Value = 19609128
$ws.Range("K136").Value = 58827384
$ws.Range("M136").Value = -58824834
